$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing rows 2..161 down to 3..162
$ws.Rows("2:2").Insert()

# New row 2 mirrors the (now shifted) old row 2 data, which landed on row 3,
# except for the date in column A which advances by one day.
$ws.Range("A2:F2").Value = $ws.Range("A3:F3").Value
$ws.Range("A2").Value = "14-01-2026"
